# Implemented excel data in Digital assessment TC
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) GradeNineVideoList: this sheet was the active tab before the edit
#    (activeTab=4, tabSelected="1", selection G6). Re-establish its
#    selection first so that it's preserved correctly once we move the
#    active tab to GradeOneVideoList further below.
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("GradeNineVideoList")
$ws9.Activate()
$ws9.Range("G6").Select()

# ---------------------------------------------------------------------
# 2) GradeOneVideoList: rewrite the lesson/subject table (rows 2-10) and
#    update the sheet view/column formatting.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("GradeOneVideoList")

# Row 2
$ws1.Cells.Item(2,1).Value = "Bible"
$ws1.Cells.Item(2,2).Value = "Lesson 1"
$ws1.Cells.Item(2,3).Value = "Reading 1"
$ws1.Cells.Item(2,4).Value = "Reading 1"
$ws1.Cells.Item(2,6).Value = 1
$ws1.Cells.Item(2,7).Value = 2

# Row 3
$ws1.Cells.Item(3,1).Value = "Classroom Routines"
$ws1.Cells.Item(3,2).Value = "Lesson 1"
$ws1.Cells.Item(3,3).Value = "Bible 1"
$ws1.Cells.Item(3,4).Value = "Bible 1"
$ws1.Cells.Item(3,6).Value = 1
$ws1.Cells.Item(3,7).Value = 2

# Row 4
$ws1.Cells.Item(4,1).Value = "Seatwork Explanation (Cursive)"
$ws1.Cells.Item(4,2).Value = "Lesson 1"
$ws1.Cells.Item(4,3).Value = "Spelling 1"
$ws1.Cells.Item(4,4).Value = "Spelling 1"
$ws1.Cells.Item(4,6).Value = 1
$ws1.Cells.Item(4,7).Value = 2

# Row 5
$ws1.Cells.Item(5,1).Value = "Phonics/Language"
$ws1.Cells.Item(5,2).Value = "Lesson 1"
$ws1.Cells.Item(5,3).Value = "Phonics 1"
$ws1.Cells.Item(5,4).Value = "Phonics 1"
$ws1.Cells.Item(5,6).Value = 1
$ws1.Cells.Item(5,7).Value = 2

# Row 6
$ws1.Cells.Item(6,1).Value = "Cursive Writing"
$ws1.Cells.Item(6,2).Value = "Lesson 1"
$ws1.Cells.Item(6,3).Value = "Activities 1"
$ws1.Cells.Item(6,4).Value = "Activities 1"
$ws1.Cells.Item(6,6).Value = 1
$ws1.Cells.Item(6,7).Value = 2

# Row 7
$ws1.Cells.Item(7,1).Value = "Spelling/Poetry"
$ws1.Cells.Item(7,2).Value = "Lesson 1"
$ws1.Cells.Item(7,3).Value = "Writing 1"
$ws1.Cells.Item(7,4).Value = "Writing 1"
$ws1.Cells.Item(7,6).Value = 1
$ws1.Cells.Item(7,7).Value = 2

# Row 8
$ws1.Cells.Item(8,1).Value = "Arithmetic"
$ws1.Cells.Item(8,2).Value = "Lesson 1"
$ws1.Cells.Item(8,3).Value = "Seatwork 1"
$ws1.Cells.Item(8,4).Value = "Seatwork 1"
$ws1.Cells.Item(8,6).Value = 1
$ws1.Cells.Item(8,7).Value = 2

# Row 9
$ws1.Cells.Item(9,1).Value = "Combination Practice"
$ws1.Cells.Item(9,2).Value = "Lesson 1"
$ws1.Cells.Item(9,3).Value = "Arithmetic 1"
$ws1.Cells.Item(9,4).Value = "Arithmetic 1"
$ws1.Cells.Item(9,6).Value = 1
$ws1.Cells.Item(9,7).Value = 2

# Row 10
$ws1.Cells.Item(10,1).Value = "Activity Time"
$ws1.Cells.Item(10,2).Value = "Lesson 1"
$ws1.Cells.Item(10,3).Value = "Classroom Routines 1"
$ws1.Cells.Item(10,4).Value = "Classroom Routines 1"
$ws1.Cells.Item(10,6).Value = 1
$ws1.Cells.Item(10,7).Value = 2

# Column D (rows 2-10) loses its left-border formatting in the new layout.
$ws1.Range("D2:D10").Borders.LineStyle = -4142

# New bestFit-style widths for columns C and D.
$ws1.Columns.Item(3).ColumnWidth = 26.83
$ws1.Columns.Item(4).ColumnWidth = 33.6

# Make GradeOneVideoList the active sheet/tab with the new selection.
$ws1.Activate()
$ws1.Range("F8").Select()
